# Generate Report for Archive
# - Flip every "Ready for handoff" status cell to "In Translation"
#   (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4)
# - Narrow the columns that held the status text to fit the new
#   (shorter) string.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # Cast to string explicitly (and keep the literal on the left) so a
        # boolean-valued cell (e.g. "True"/"False" columns) never gets
        # coerced into matching a non-empty comparison string.
        if ($oldStatus -eq [string]$cell.Text) {
            $cell.Value = $newStatus
        }
    }
}

# The "zh-cn"/"de-de" Status columns (and the Overview sheet's mirrored
# zh-cn/de-de columns) were sized for the old, longer "Ready for handoff"
# text. Re-fit them now that the text is shorter.
$narrowWidth = 12.5

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $narrowWidth
$wsOverview.Columns.Item(6).ColumnWidth = $narrowWidth

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = $narrowWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = $narrowWidth
